$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's format (bold, centered, bordered) onto the new H1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
